$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental flag: false -> true
$ws.Range("B7").Value = "true"

# Date updated
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"

# Case Sensitive: was blank, now "false"
$ws.Range("B14").Value = "false"
